$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2: A2 = 1 (was 0), B2 stays 209
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 209

# Update row 3: A3 = 0 (was 1), B3 = 89 (was 66)
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 89

# Remove row 4 entirely (was A4=2, B4=23)
$ws.Range("A4:B4").Value = $null
$ws.Rows.Item(4).Delete()
